# Added RST+TX+RX to Expansion Header
# The expansion header (J4) is upgraded from a 1x04 connector to a 1x07
# connector (Samtec SSW-107-01-T-S) to carry RST, TXD and RXD in addition
# to the existing signals. Two extra 1K resistors (R31, R32) are fitted
# so the FT232RL's TXD/RXD lines can be safely overridden by an expansion
# device, so the R18, R25-R26 resistor line is updated to reflect the
# extra references and increased quantity.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 25: J4 expansion header connector, swap 1x04 -> 1x07 part ---
$ws.Range("C25").Value = "SSW-107-01-T-S"
$ws.Range("F25").Value = 0.735
$ws.Range("H25").Value = "200-SSW10701TS"
$ws.Range("I25").Value = "2667434"
$ws.Range("L25").Value = "1x07 2.54mm Square Header Recepticle, Vertical, THT"
$ws.Range("K25").Value = "Connector_PinSocket_2.54mm:PinSocket_1x07_P2.54mm_Vertical"

# --- Row 34: 1kR resistors, now also used for R31-R32, qty 3 -> 5 ---
$ws.Range("A34").Value = "R18, R25-R26, R31-R32"
$ws.Range("D34").Value = 5

# Update the active selection to A35, matching the final cursor position
# left by the author after editing the resistor row.
$ws.Range("A35").Select()

$wb.Save()
